$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "degHh927"
$ws.Range("B2").Value = 23080135
$ws.Range("C2").Value = "rchbrlj12"
$ws.Range("D2").Value = "h%8eB5!N"
$ws.Range("F2").Value = "BDQEAMiI"
$ws.Range("G2").Value = "sYkg"
